$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")
Write-Host $ws.Name
